$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the date values for rows 2-4 in column A
$ws.Range("A2").Value = "3/21/2025"
$ws.Range("A3").Value = "3/26/2025"
$ws.Range("A4").Value = "4/2/2025"

# Update the selected cell/range on the sheet (as recorded in the sheetView)
$ws.Range("E8").Select()
